$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# --- New data for rows 2-17 (row 1 is the header and is unchanged) ---
$rows = @(
    @{ A="2026-02-03 06:47:40"; B="Gmail、スプレッドシート、Google Driveを連携した 業務効率化システム開発"; C="システム開発"; D="100,000 円 ~ 200,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5484339"; G=453; H="🔥AI,Ai ◆効率化,開発" }
    @{ A="2026-02-03 06:47:40"; B="建設・土木業界向け施工機械のAI自動制御・アタッチメント開発を支援してくださるエンジニア募集"; C="システム開発"; D="200,000 円 ~ 300,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5434128"; G=368; H="🔥AI,Ai ◆開発" }
    @{ A="2026-02-03 06:47:40"; B="大手製造業向け センサー画像解析・高画質化のR&Dを支援するAIエンジニア募集(画像生成/超解像)"; C="システム開発"; D="300,000 円 ~ 500,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5427956"; G=310; H="🔥AI,Ai" }
    @{ A="2026-02-03 06:47:40"; B="【急募】製造業向け「製造副産物」の状態(硬度)判定AIのフィジビリティ検証(画像認識/動画解析)"; C="システム開発"; D="200,000 円 ~ 300,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5439158"; G=303; H="🔥AI,Ai" }
    @{ A="2026-02-03 06:47:40"; B="急募 【急募】日繰資金繰表自動算出ツールの開発依頼"; C="システム開発"; D="20,000 円 ~ 50,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5484458"; G=123; H="◆ツール,開発" }
    @{ A="2026-02-03 06:47:40"; B="(仕様削減)【受注メールを元にしたスクレピング&抽出情報管理ツール】"; C="システム開発"; D="50,000 円 ~ 100,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5484248"; G=98; H="◆ツール ◇管理" }
    @{ A="2026-02-03 06:47:40"; B="ストレスチェック集団分析自動化・レポート出力システムを探しています"; C="システム開発"; D="20,000 円 ~ 50,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5484413"; G=98; H="◆自動化" }
    @{ A="2026-02-03 06:47:40"; B="M5Stackを用いたAC100V電気ケトルの制御回路・試作機開発"; C="システム開発"; D="50,000 円 ~ 100,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5484617"; G=68; H="◆開発" }
    @{ A="2026-02-03 06:47:40"; B="リアルタイム顔変換開発案件 技術判断を整理してくれる方(短期スポット)募集"; C="システム開発"; D="100,000 円 ~ 200,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5484376"; G=68; H="◆開発" }
    @{ A="2026-02-03 06:47:40"; B="【店舗DX】LINE自動予約システム導入で集客力向上を目指す"; C="システム開発"; D="20,000 円 ~ 50,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5484779"; G=28; H=$null }
    @{ A="2026-02-03 06:47:40"; B="【急募】シェアポイント上のリストから自動での各種文書作成のプロフェッショナルを探しています!"; C="システム開発"; D="200,000 円 ~ 300,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5484682"; G=18; H=$null }
    @{ A="2026-02-03 06:47:40"; B="【急募】Unity/Photonでのメタバースマルチプレイヤー機能実装"; C="システム開発"; D="200,000 円 ~ 300,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5484436"; G=18; H=$null }
    @{ A="2026-02-03 06:47:40"; B="防災プラットフォームの作成"; C="システム開発"; D="100,000 円 ~ 200,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5484389"; G=18; H=$null }
    @{ A="2026-02-03 06:47:40"; B="限定公開 限定公開の仕事"; C="システム開発"; D="50,000 円 ~ 100,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5484369"; G=18; H=$null }
    @{ A="2026-02-03 06:47:40"; B="【フリーランス必見】エンジニア支援サービスのご紹介!"; C="システム開発"; D="10,000 円 ~ 20,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5484757"; G=10; H=$null }
    @{ A="2026-02-03 06:47:40"; B="Claude Codeにプロンプトを入れてください"; C="システム開発"; D="1,000 ~ 5,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5484754"; G=10; H=$null }
)

# Remove the now-stale hyperlink collection up front; it will be rebuilt
# below once the final row count/URLs are known.
$ws.Hyperlinks.Delete()

# Delete the 4 trailing rows (old rows 18-21) - the newest scrape only has
# 16 job postings (rows 2-17), four fewer than before.
$ws.Rows.Item(21).Delete()
$ws.Rows.Item(20).Delete()
$ws.Rows.Item(19).Delete()
$ws.Rows.Item(18).Delete()

# Write the refreshed rows.
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data.A
    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 3).Value = $data.C
    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 5).Value = $data.E
    $ws.Cells.Item($r, 6).Value = $data.F
    $ws.Cells.Item($r, 7).Value = $data.G
    if ($data.H -eq $null) {
        $ws.Cells.Item($r, 8).ClearContents()
    } else {
        $ws.Cells.Item($r, 8).Value = $data.H
    }
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $data.F)
}

# Column-width tweaks (D: 30 -> 28, H: 19 -> 16). ColumnWidth uses a
# chars-to-stored-width offset of +5/6, so subtract that to land on the
# exact target stored width.
$ws.Columns.Item(4).ColumnWidth = 28 - 5/6
$ws.Columns.Item(8).ColumnWidth = 16 - 5/6
